$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 8 (Magnesium chloride unit price): update loading/lower/upper values
# and convert the previously-formula-driven Lower/Upper cells into plain literal values.
$ws.Range("E8").Value = 0.38
$ws.Range("G8").Value = 0.349
$ws.Range("I8").Value = 0.411

# Row 9 (Zinc sulfate unit price): same treatment
$ws.Range("E9").Value = 0.795
$ws.Range("G9").Value = 0.657
$ws.Range("I9").Value = 0.931

# The helper "match" column (Q) for row 8 used to carry its own one-off
# shared formula; re-entering row 8 consolidates it back into the common
# Q5:Q32 shared-formula group (same computed result, 1).
$ws.Range("Q8").Value = 1

$ws.Rows("8:9").Select()
